$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "tahun" -> "tahun_anggaran" (column S, header row 1)
$ws.Range("S1").Value = "tahun_anggaran"

# Format the "nik" column (G) as Text so long ID numbers aren't mangled
$ws.Range("G1").NumberFormat = "@"

# New column S width (for the imported disability data) and view position
$ws.Columns("S").ColumnWidth = 17.4444444444444

$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("R20").Select()
